# "clean up before submitting"
#
# The second raw data column (old column C, "min78 m2") is removed; the
# existing km2 column (old column D, =B/1e6) slides left into column C and
# every downstream helper formula / the line chart follow it automatically.
# The underlying raw measurements in column B are also refreshed with a
# newer run of the analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Refresh the raw data in column B (rows 2-24) with the re-run numbers.
# ---------------------------------------------------------------------------
$newB = @(
    2982275336.8035598,
    3342535226.2813601,
    5276671735.4082403,
    3912847746.0925002,
    2614731722.5745702,
    3032727768.3896298,
    2486222167.1809502,
    6131422841.9177303,
    5975433671.36724,
    10485177827.728201,
    8956578656.1823006,
    7212638739.92208,
    11804068916.794701,
    6444022418.6795197,
    10446949205.1812,
    7006335060.1411695,
    10029907654.908501,
    3183702802.8784299,
    10196072328.743099,
    8495792917.0384302,
    6015131586.9759197,
    11034799341.0753,
    6988325407.7599096
)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $newB[$i]
}

# ---------------------------------------------------------------------------
# 2. Before removing the column, work out the chart's new on-sheet position
#    (it gets moved/shrunk once the stray data column disappears). We
#    compute this from the *current* (pre-delete) geometry, since the target
#    anchor cells/offsets were authored against that same column layout.
# ---------------------------------------------------------------------------
function ConvertEmuToPoints($emu) {
    return $emu / 12700.0
}

# twoCellAnchor is 0-indexed; Cells.Item is 1-indexed.
$fromCell = $ws.Cells.Item(16, 6)    # row15 -> 16, col5 -> 6 (column F)
$toCell = $ws.Cells.Item(31, 20)     # row30 -> 31, col19 -> 20 (column T)

$newLeft = $fromCell.Left + (ConvertEmuToPoints 497204)
$newTop = $fromCell.Top + (ConvertEmuToPoints 114300)
$newRight = $toCell.Left + (ConvertEmuToPoints 295274)
$newBottom = $toCell.Top + (ConvertEmuToPoints 9525)

# ---------------------------------------------------------------------------
# 3. Delete the obsolete raw-data column (old column C, "min78 m2").
#    This shifts: D (km2) -> C, G/H/I (helper labels/averages) -> F/G/H,
#    and Excel keeps every formula referring to the shifted cells correct.
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Delete()

# ---------------------------------------------------------------------------
# 4. The third summary ratio row (old "=(H4-H2)/H2" in I4, now H4) only
#    existed for the last bucket. Add the matching ratio for the middle
#    bucket (now row 3), matching H4's percentage formatting.
# ---------------------------------------------------------------------------
$ws.Range("H3").Formula = "=(G3-G2)/G2"
$ws.Range("H3").NumberFormat = $ws.Range("H4").NumberFormat

# ---------------------------------------------------------------------------
# 5. Re-point the chart series at the (now-shifted) km2 column and move /
#    resize the chart to its new location.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$A`$24,Sheet1!`$C`$2:`$C`$24,1)"

$co.Left = $newLeft
$co.Top = $newTop
$co.Width = $newRight - $newLeft
$co.Height = $newBottom - $newTop

# ---------------------------------------------------------------------------
# 6. Update the active selection to match what was left selected.
# ---------------------------------------------------------------------------
$ws.Range("D11").Select()
